# Applies the numeric odds updates for the 2026-01-07 Betfair Back/Lay workbook.
# Each assignment below corresponds to one changed cell (row, column) in the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 10).Value = 5.4  # J2: 5.3 -> 5.4
$ws.Cells.Item(2, 16).Value = 1.86  # P2: 1.85 -> 1.86
$ws.Cells.Item(2, 21).Value = 1.6  # U2: 1.58 -> 1.6
$ws.Cells.Item(2, 23).Value = 3.65  # W2: 3.6 -> 3.65
$ws.Cells.Item(2, 27).Value = 740  # AA2: 750 -> 740
$ws.Cells.Item(2, 41).Value = 610  # AO2: 620 -> 610
# Row 3
$ws.Cells.Item(3, 17).Value = 2.04  # Q3: 2.02 -> 2.04
# Row 4
$ws.Cells.Item(4, 8).Value = 8.4  # H4: 8.199999999999999 -> 8.4
$ws.Cells.Item(4, 9).Value = 8.6  # I4: 8.4 -> 8.6
$ws.Cells.Item(4, 12).Value = 1.22  # L4: 0 -> 1.22
$ws.Cells.Item(4, 13).Value = 1.03  # M4: 1.02 -> 1.03
$ws.Cells.Item(4, 17).Value = 1.47  # Q4: 1.46 -> 1.47
$ws.Cells.Item(4, 18).Value = 1.82  # R4: 1.83 -> 1.82
$ws.Cells.Item(4, 21).Value = 2.4  # U4: 2.38 -> 2.4
$ws.Cells.Item(4, 22).Value = 1.13  # V4: 0 -> 1.13
$ws.Cells.Item(4, 23).Value = 3.35  # W4: 0 -> 3.35
$ws.Cells.Item(4, 24).Value = 32  # X4: 36 -> 32
$ws.Cells.Item(4, 26).Value = 80  # Z4: 85 -> 80
$ws.Cells.Item(4, 27).Value = 220  # AA4: 260 -> 220
$ws.Cells.Item(4, 30).Value = 29  # AD4: 32 -> 29
$ws.Cells.Item(4, 31).Value = 90  # AE4: 95 -> 90
$ws.Cells.Item(4, 33).Value = 10  # AG4: 10.5 -> 10
$ws.Cells.Item(4, 35).Value = 75  # AI4: 80 -> 75
$ws.Cells.Item(4, 38).Value = 24  # AL4: 25 -> 24
$ws.Cells.Item(4, 39).Value = 80  # AM4: 85 -> 80
$ws.Cells.Item(4, 41).Value = 75  # AO4: 85 -> 75
# Row 5
$ws.Cells.Item(5, 10).Value = 3.8  # J5: 3.85 -> 3.8
$ws.Cells.Item(5, 12).Value = 1.31  # L5: 0 -> 1.31
$ws.Cells.Item(5, 16).Value = 2.24  # P5: 2.22 -> 2.24
$ws.Cells.Item(5, 18).Value = 1.48  # R5: 1.47 -> 1.48
$ws.Cells.Item(5, 21).Value = 2.42  # U5: 2.38 -> 2.42
$ws.Cells.Item(5, 22).Value = 1.39  # V5: 0 -> 1.39
$ws.Cells.Item(5, 23).Value = 1.83  # W5: 0 -> 1.83
$ws.Cells.Item(5, 24).Value = 17.5  # X5: 19 -> 17.5
$ws.Cells.Item(5, 27).Value = 65  # AA5: 75 -> 65
$ws.Cells.Item(5, 30).Value = 14  # AD5: 15 -> 14
$ws.Cells.Item(5, 31).Value = 36  # AE5: 38 -> 36
$ws.Cells.Item(5, 32).Value = 14.5  # AF5: 15.5 -> 14.5
$ws.Cells.Item(5, 37).Value = 21  # AK5: 22 -> 21
$ws.Cells.Item(5, 38).Value = 30  # AL5: 32 -> 30
$ws.Cells.Item(5, 41).Value = 29  # AO5: 30 -> 29
# Row 6
$ws.Cells.Item(6, 12).Value = 1.43  # L6: 0 -> 1.43
$ws.Cells.Item(6, 13).Value = 1.09  # M6: 1.08 -> 1.09
$ws.Cells.Item(6, 17).Value = 2.2  # Q6: 2.18 -> 2.2
$ws.Cells.Item(6, 22).Value = 1.21  # V6: 0 -> 1.21
$ws.Cells.Item(6, 23).Value = 2.22  # W6: 0 -> 2.22
$ws.Cells.Item(6, 24).Value = 12  # X6: 12.5 -> 12
$ws.Cells.Item(6, 26).Value = 40  # Z6: 42 -> 40
$ws.Cells.Item(6, 27).Value = 150  # AA6: 190 -> 150
$ws.Cells.Item(6, 30).Value = 21  # AD6: 22 -> 21
$ws.Cells.Item(6, 32).Value = 9.4  # AF6: 10.5 -> 9.4
$ws.Cells.Item(6, 33).Value = 10  # AG6: 10.5 -> 10
$ws.Cells.Item(6, 36).Value = 18  # AJ6: 19.5 -> 18
$ws.Cells.Item(6, 39).Value = 160  # AM6: 180 -> 160
$ws.Cells.Item(6, 40).Value = 14  # AN6: 14.5 -> 14
$ws.Cells.Item(6, 41).Value = 120  # AO6: 150 -> 120
# Row 7
$ws.Cells.Item(7, 6).Value = 3.3  # F7: 3.25 -> 3.3
$ws.Cells.Item(7, 7).Value = 3.35  # G7: 3.3 -> 3.35
$ws.Cells.Item(7, 8).Value = 2.44  # H7: 2.42 -> 2.44
$ws.Cells.Item(7, 9).Value = 2.46  # I7: 2.44 -> 2.46
$ws.Cells.Item(7, 12).Value = 1.39  # L7: 0 -> 1.39
$ws.Cells.Item(7, 14).Value = 3.85  # N7: 3.9 -> 3.85
$ws.Cells.Item(7, 22).Value = 1.68  # V7: 0 -> 1.68
$ws.Cells.Item(7, 23).Value = 1.42  # W7: 0 -> 1.42
$ws.Cells.Item(7, 24).Value = 13.5  # X7: 14 -> 13.5
$ws.Cells.Item(7, 25).Value = 10.5  # Y7: 11 -> 10.5
$ws.Cells.Item(7, 26).Value = 15  # Z7: 16 -> 15
$ws.Cells.Item(7, 27).Value = 32  # AA7: 34 -> 32
$ws.Cells.Item(7, 28).Value = 13  # AB7: 13.5 -> 13
$ws.Cells.Item(7, 30).Value = 11  # AD7: 11.5 -> 11
$ws.Cells.Item(7, 32).Value = 22  # AF7: 23 -> 22
$ws.Cells.Item(7, 33).Value = 13.5  # AG7: 14.5 -> 13.5
$ws.Cells.Item(7, 35).Value = 38  # AI7: 40 -> 38
$ws.Cells.Item(7, 38).Value = 46  # AL7: 48 -> 46
$ws.Cells.Item(7, 39).Value = 90  # AM7: 1000 -> 90
$ws.Cells.Item(7, 41).Value = 20  # AO7: 21 -> 20
# Row 8
$ws.Cells.Item(8, 6).Value = 1.92  # F8: 1.91 -> 1.92
$ws.Cells.Item(8, 7).Value = 1.94  # G8: 1.92 -> 1.94
$ws.Cells.Item(8, 12).Value = 1.41  # L8: 0 -> 1.41
$ws.Cells.Item(8, 19).Value = 3.9  # S8: 3.85 -> 3.9
$ws.Cells.Item(8, 22).Value = 1.27  # V8: 0 -> 1.27
$ws.Cells.Item(8, 23).Value = 2.06  # W8: 0 -> 2.06
$ws.Cells.Item(8, 24).Value = 13  # X8: 12.5 -> 13
$ws.Cells.Item(8, 25).Value = 15.5  # Y8: 15 -> 15.5
$ws.Cells.Item(8, 26).Value = 32  # Z8: 34 -> 32
$ws.Cells.Item(8, 27).Value = 120  # AA8: 130 -> 120
$ws.Cells.Item(8, 28).Value = 8.199999999999999  # AB8: 8.4 -> 8.199999999999999
$ws.Cells.Item(8, 30).Value = 18  # AD8: 18.5 -> 18
$ws.Cells.Item(8, 32).Value = 11  # AF8: 11.5 -> 11
$ws.Cells.Item(8, 34).Value = 21  # AH8: 22 -> 21
$ws.Cells.Item(8, 35).Value = 75  # AI8: 70 -> 75
$ws.Cells.Item(8, 36).Value = 21  # AJ8: 22 -> 21
$ws.Cells.Item(8, 37).Value = 21  # AK8: 20 -> 21
$ws.Cells.Item(8, 39).Value = 110  # AM8: 100 -> 110
$ws.Cells.Item(8, 40).Value = 15  # AN8: 15.5 -> 15
$ws.Cells.Item(8, 41).Value = 75  # AO8: 100 -> 75
# Row 9
$ws.Cells.Item(9, 6).Value = 3.65  # F9: 3.6 -> 3.65
$ws.Cells.Item(9, 7).Value = 3.7  # G9: 3.65 -> 3.7
$ws.Cells.Item(9, 8).Value = 2.12  # H9: 2.16 -> 2.12
$ws.Cells.Item(9, 9).Value = 2.14  # I9: 2.18 -> 2.14
$ws.Cells.Item(9, 10).Value = 3.8  # J9: 3.75 -> 3.8
$ws.Cells.Item(9, 11).Value = 3.85  # K9: 3.8 -> 3.85
$ws.Cells.Item(9, 12).Value = 1.34  # L9: 0 -> 1.34
$ws.Cells.Item(9, 20).Value = 1.73  # T9: 1.72 -> 1.73
$ws.Cells.Item(9, 22).Value = 1.87  # V9: 0 -> 1.87
$ws.Cells.Item(9, 23).Value = 1.37  # W9: 0 -> 1.37
$ws.Cells.Item(9, 24).Value = 16  # X9: 17 -> 16
$ws.Cells.Item(9, 26).Value = 13.5  # Z9: 14.5 -> 13.5
$ws.Cells.Item(9, 27).Value = 26  # AA9: 27 -> 26
$ws.Cells.Item(9, 28).Value = 15  # AB9: 15.5 -> 15
$ws.Cells.Item(9, 29).Value = 8.199999999999999  # AC9: 8.4 -> 8.199999999999999
$ws.Cells.Item(9, 31).Value = 21  # AE9: 20 -> 21
$ws.Cells.Item(9, 32).Value = 26  # AF9: 27 -> 26
$ws.Cells.Item(9, 33).Value = 14.5  # AG9: 15 -> 14.5
$ws.Cells.Item(9, 39).Value = 80  # AM9: 75 -> 80
$ws.Cells.Item(9, 41).Value = 14  # AO9: 15 -> 14
# Row 10
$ws.Cells.Item(10, 6).Value = 2.54  # F10: 2.52 -> 2.54
$ws.Cells.Item(10, 7).Value = 2.56  # G10: 2.54 -> 2.56
$ws.Cells.Item(10, 12).Value = 1.47  # L10: 0 -> 1.47
$ws.Cells.Item(10, 22).Value = 1.41  # V10: 0 -> 1.41
$ws.Cells.Item(10, 23).Value = 1.64  # W10: 0 -> 1.64
$ws.Cells.Item(10, 24).Value = 9.800000000000001  # X10: 10.5 -> 9.800000000000001
$ws.Cells.Item(10, 26).Value = 21  # Z10: 22 -> 21
$ws.Cells.Item(10, 27).Value = 65  # AA10: 70 -> 65
$ws.Cells.Item(10, 30).Value = 14.5  # AD10: 15 -> 14.5
$ws.Cells.Item(10, 31).Value = 46  # AE10: 48 -> 46
$ws.Cells.Item(10, 32).Value = 14.5  # AF10: 15 -> 14.5
$ws.Cells.Item(10, 33).Value = 12  # AG10: 12.5 -> 12
$ws.Cells.Item(10, 35).Value = 65  # AI10: 70 -> 65
$ws.Cells.Item(10, 39).Value = 140  # AM10: 170 -> 140
$ws.Cells.Item(10, 40).Value = 30  # AN10: 32 -> 30
$ws.Cells.Item(10, 41).Value = 55  # AO10: 60 -> 55
# Row 11
$ws.Cells.Item(11, 12).Value = 1.5  # L11: 0 -> 1.5
$ws.Cells.Item(11, 22).Value = 1.37  # V11: 0 -> 1.37
$ws.Cells.Item(11, 23).Value = 1.71  # W11: 0 -> 1.71
$ws.Cells.Item(11, 24).Value = 9.4  # X11: 10 -> 9.4
$ws.Cells.Item(11, 25).Value = 11  # Y11: 11.5 -> 11
$ws.Cells.Item(11, 27).Value = 70  # AA11: 75 -> 70
$ws.Cells.Item(11, 30).Value = 15.5  # AD11: 16 -> 15.5
$ws.Cells.Item(11, 31).Value = 50  # AE11: 55 -> 50
$ws.Cells.Item(11, 32).Value = 14  # AF11: 14.5 -> 14
$ws.Cells.Item(11, 33).Value = 11.5  # AG11: 12 -> 11.5
$ws.Cells.Item(11, 36).Value = 32  # AJ11: 36 -> 32
$ws.Cells.Item(11, 37).Value = 30  # AK11: 32 -> 30
$ws.Cells.Item(11, 39).Value = 160  # AM11: 170 -> 160
$ws.Cells.Item(11, 40).Value = 29  # AN11: 30 -> 29
# Row 12
$ws.Cells.Item(12, 6).Value = 9.6  # F12: 9.199999999999999 -> 9.6
$ws.Cells.Item(12, 7).Value = 9.800000000000001  # G12: 9.6 -> 9.800000000000001
$ws.Cells.Item(12, 10).Value = 5.5  # J12: 5.4 -> 5.5
$ws.Cells.Item(12, 11).Value = 5.6  # K12: 5.5 -> 5.6
$ws.Cells.Item(12, 12).Value = 1.29  # L12: 0 -> 1.29
$ws.Cells.Item(12, 14).Value = 5  # N12: 5.1 -> 5
$ws.Cells.Item(12, 15).Value = 1.23  # O12: 1.22 -> 1.23
$ws.Cells.Item(12, 16).Value = 2.36  # P12: 2.38 -> 2.36
$ws.Cells.Item(12, 20).Value = 1.98  # T12: 1.97 -> 1.98
$ws.Cells.Item(12, 22).Value = 3.45  # V12: 0 -> 3.45
$ws.Cells.Item(12, 23).Value = 1.11  # W12: 0 -> 1.11
$ws.Cells.Item(12, 24).Value = 21  # X12: 23 -> 21
$ws.Cells.Item(12, 26).Value = 8.6  # Z12: 8.800000000000001 -> 8.6
$ws.Cells.Item(12, 27).Value = 11.5  # AA12: 12 -> 11.5
$ws.Cells.Item(12, 30).Value = 9.800000000000001  # AD12: 10.5 -> 9.800000000000001
$ws.Cells.Item(12, 31).Value = 14  # AE12: 15 -> 14
$ws.Cells.Item(12, 32).Value = 85  # AF12: 1000 -> 85
$ws.Cells.Item(12, 35).Value = 34  # AI12: 36 -> 34
$ws.Cells.Item(12, 36).Value = 370  # AJ12: 380 -> 370
$ws.Cells.Item(12, 37).Value = 160  # AK12: 170 -> 160
$ws.Cells.Item(12, 38).Value = 120  # AL12: 1000 -> 120
$ws.Cells.Item(12, 39).Value = 150  # AM12: 160 -> 150
$ws.Cells.Item(12, 40).Value = 170  # AN12: 210 -> 170
# Row 13
$ws.Cells.Item(13, 6).Value = 1.73  # F13: 1.72 -> 1.73
$ws.Cells.Item(13, 7).Value = 1.74  # G13: 1.73 -> 1.74
$ws.Cells.Item(13, 9).Value = 5.8  # I13: 5.7 -> 5.8
$ws.Cells.Item(13, 12).Value = 1.36  # L13: 0 -> 1.36
$ws.Cells.Item(13, 16).Value = 2.12  # P13: 2.14 -> 2.12
$ws.Cells.Item(13, 19).Value = 3.15  # S13: 3.2 -> 3.15
$ws.Cells.Item(13, 20).Value = 1.84  # T13: 1.85 -> 1.84
$ws.Cells.Item(13, 22).Value = 1.21  # V13: 0 -> 1.21
$ws.Cells.Item(13, 23).Value = 2.34  # W13: 0 -> 2.34
$ws.Cells.Item(13, 24).Value = 16  # X13: 17 -> 16
$ws.Cells.Item(13, 25).Value = 20  # Y13: 21 -> 20
$ws.Cells.Item(13, 26).Value = 42  # Z13: 44 -> 42
$ws.Cells.Item(13, 27).Value = 140  # AA13: 170 -> 140
$ws.Cells.Item(13, 29).Value = 9  # AC13: 8.800000000000001 -> 9
$ws.Cells.Item(13, 30).Value = 21  # AD13: 22 -> 21
$ws.Cells.Item(13, 31).Value = 70  # AE13: 75 -> 70
$ws.Cells.Item(13, 32).Value = 10.5  # AF13: 11 -> 10.5
$ws.Cells.Item(13, 33).Value = 9.6  # AG13: 9.800000000000001 -> 9.6
$ws.Cells.Item(13, 35).Value = 70  # AI13: 90 -> 70
$ws.Cells.Item(13, 36).Value = 17  # AJ13: 17.5 -> 17
$ws.Cells.Item(13, 37).Value = 16.5  # AK13: 17.5 -> 16.5
$ws.Cells.Item(13, 38).Value = 32  # AL13: 34 -> 32
$ws.Cells.Item(13, 39).Value = 100  # AM13: 110 -> 100
$ws.Cells.Item(13, 40).Value = 9.6  # AN13: 9.800000000000001 -> 9.6
$ws.Cells.Item(13, 41).Value = 80  # AO13: 100 -> 80
# Row 14
$ws.Cells.Item(14, 6).Value = 5.3  # F14: 5.4 -> 5.3
$ws.Cells.Item(14, 7).Value = 5.4  # G14: 5.5 -> 5.4
$ws.Cells.Item(14, 12).Value = 1.33  # L14: 0 -> 1.33
$ws.Cells.Item(14, 13).Value = 1.06  # M14: 1.05 -> 1.06
$ws.Cells.Item(14, 14).Value = 4.5  # N14: 4.4 -> 4.5
$ws.Cells.Item(14, 17).Value = 1.82  # Q14: 1.83 -> 1.82
$ws.Cells.Item(14, 18).Value = 1.47  # R14: 1.46 -> 1.47
$ws.Cells.Item(14, 19).Value = 3  # S14: 3.05 -> 3
$ws.Cells.Item(14, 20).Value = 1.8  # T14: 1.79 -> 1.8
$ws.Cells.Item(14, 22).Value = 2.3  # V14: 0 -> 2.3
$ws.Cells.Item(14, 23).Value = 1.22  # W14: 0 -> 1.22
$ws.Cells.Item(14, 24).Value = 17.5  # X14: 19 -> 17.5
$ws.Cells.Item(14, 27).Value = 17.5  # AA14: 18 -> 17.5
$ws.Cells.Item(14, 30).Value = 9.6  # AD14: 9.800000000000001 -> 9.6
$ws.Cells.Item(14, 31).Value = 16.5  # AE14: 17.5 -> 16.5
$ws.Cells.Item(14, 33).Value = 20  # AG14: 21 -> 20
$ws.Cells.Item(14, 35).Value = 32  # AI14: 34 -> 32
$ws.Cells.Item(14, 36).Value = 130  # AJ14: 150 -> 130
$ws.Cells.Item(14, 37).Value = 65  # AK14: 70 -> 65
$ws.Cells.Item(14, 38).Value = 65  # AL14: 70 -> 65
$ws.Cells.Item(14, 39).Value = 95  # AM14: 90 -> 95
$ws.Cells.Item(14, 40).Value = 65  # AN14: 75 -> 65
